$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data for this run.

# Row 2
$ws.Range("D2").Value = "42.451.30"
$ws.Range("E2").Value = "  +1.29%  "

# Row 3
$ws.Range("D3").Value = "2.299.04"
$ws.Range("E3").Value = "  +0.93%  "

# Row 4
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.63"
$ws.Range("E5").Value = "  +1.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.69"
$ws.Range("E6").Value = "  -2.69%  "

# Row 7
$ws.Range("E7").Value = "  +0.67%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  +0.31%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.59"
$ws.Range("E10").Value = "  -1.60%  "

# Row 11
$ws.Range("E11").Value = "  +0.07%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.40"
$ws.Range("E12").Value = "  +1.80%  "

# Row 13
$ws.Range("E13").Value = "  +0.78%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.961"
$ws.Range("E14").Value = "  +0.00%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.23"
$ws.Range("E15").Value = "  -1.11%  "

# Row 16
$ws.Range("D16").Value = "2.647.52"
$ws.Range("E16").Value = "  +0.89%  "

# Row 17
$ws.Range("D17").Value = "2.304.65"
$ws.Range("E17").Value = "  +1.82%  "

# Row 18
$ws.Range("D18").Value = "42.398.91"
$ws.Range("E18").Value = "  +1.33%  "

# Row 19
$ws.Range("E19").Value = "  -1.58%  "

# Row 20
$ws.Range("E20").Value = "  +1.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.32"
$ws.Range("E21").Value = "  +0.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.64"
$ws.Range("E22").Value = "  +25.63%  "

# Row 23
$ws.Range("E23").Value = "  +3.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "274.71"
$ws.Range("E24").Value = "  +7.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.26"
$ws.Range("E25").Value = "  -1.25%  "

# Row 26
$ws.Range("E26").Value = "  -0.40%  "

# Row 27
$ws.Range("E27").Value = "  -0.73%  "

# Row 28
$ws.Range("E28").Value = "  +3.58%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.70"
$ws.Range("E29").Value = "  +0.72%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.43"
$ws.Range("E30").Value = "  +5.84%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.68"
$ws.Range("E31").Value = "  +0.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0873"
$ws.Range("E32").Value = "  -1.33%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.93"
$ws.Range("E33").Value = "  +3.49%  "

# Row 34
$ws.Range("E34").Value = "  +3.88%  "

# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.63"
$ws.Range("E35").Value = "  -9.75%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  +1.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.58"
$ws.Range("E37").Value = "  +1.04%  "

# Row 38
$ws.Range("E38").Value = "  +3.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.73"
$ws.Range("E39").Value = "  +3.26%  "

# Row 40
$ws.Range("E40").Value = "  +0.34%  "

# Row 41
$ws.Range("E41").Value = "  +3.24%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.99"
$ws.Range("E42").Value = "  -1.94%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "94.63"
$ws.Range("E43").Value = "  -2.48%  "

# Row 44
$ws.Range("E44").Value = "  +0.12%  "

# Row 45
$ws.Range("E45").Value = "  +0.24%  "

# Row 46
$ws.Range("E46").Value = "  -1.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "80.69"
$ws.Range("E47").Value = "  +8.10%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.46"
$ws.Range("E48").Value = "  +0.95%  "

# Row 49
$ws.Range("E49").Value = "  +0.03%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.25"
$ws.Range("E50").Value = "  -0.44%  "

# Row 51
$ws.Range("D51").Value = "1.588.76"
$ws.Range("E51").Value = "  +2.23%  "
